$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Overwrite row 16 data with the record previously in row 19 (Rosa Angelica)
$ws.Range("C16").Value = "1143394352"
$ws.Range("D16").Value = "ROSA ANGELICA BORJAS CARDENAS"
$ws.Range("E16").Value = "2403"
$ws.Range("F16").Value = 26000
$ws.Range("G16").Value = 1350000

# Remove the now-duplicate rows (old rows 17, 18, 19)
$ws.Range("17:19").EntireRow.Delete()

# Update the summary counters
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1
$ws.Range("E11").Value = 26000
